# Fruta / hortaliza, semanal
# Inserts 4 new weekly-report rows (207-210) for "Early Majestic" / "Florida King"
# Durazno lots and shifts the previously-existing rows 207-225 down to 211-229,
# matching the new sheet dimension A1:T229.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 207:225 down to 211:229 by inserting 4 blank rows at 207.
$ws.Rows("207:210").Insert()

# Static (non-varying) column values shared by every row in this data block.
$marketId = 4
$market   = "Feria Lagunitas de Puerto Montt"
$region   = "Los Lagos"
$codreg   = 10
$tipo     = "Fruta"
$prodId   = 100103
$prodName = "Frutos de hueso (carozo)"
$catId    = 100103004
$catName  = "Durazno"

function Set-DuraznoRow {
    param(
        [int]$Row,
        [double]$Fecha,
        [string]$Variedad,
        [string]$Calidad,
        [double]$Volumen,
        [double]$PrecioMin,
        [double]$PrecioMax,
        [double]$PrecioProm,
        [string]$Unidad,
        [string]$Origen,
        [double]$PrecioKg,
        [double]$KgUnidad
    )

    $ws.Cells.Item($Row, 1).Value  = $marketId
    $ws.Cells.Item($Row, 2).Value  = $market
    $ws.Cells.Item($Row, 3).Value  = $region
    $ws.Cells.Item($Row, 4).Value  = $Fecha
    $ws.Cells.Item($Row, 5).Value  = $codreg
    $ws.Cells.Item($Row, 6).Value  = $tipo
    $ws.Cells.Item($Row, 7).Value  = $prodId
    $ws.Cells.Item($Row, 8).Value  = $prodName
    $ws.Cells.Item($Row, 9).Value  = $catId
    $ws.Cells.Item($Row, 10).Value = $catName
    $ws.Cells.Item($Row, 11).Value = $Variedad
    $ws.Cells.Item($Row, 12).Value = $Calidad
    $ws.Cells.Item($Row, 13).Value = $Volumen
    $ws.Cells.Item($Row, 14).Value = $PrecioMin
    $ws.Cells.Item($Row, 15).Value = $PrecioMax
    $ws.Cells.Item($Row, 16).Value = $PrecioProm
    $ws.Cells.Item($Row, 17).Value = $Unidad
    $ws.Cells.Item($Row, 18).Value = $Origen
    $ws.Cells.Item($Row, 19).Value = $PrecioKg
    $ws.Cells.Item($Row, 20).Value = $KgUnidad
}

# Set-DuraznoRow Row Fecha Variedad Calidad Volumen PrecioMin PrecioMax PrecioProm Unidad Origen PrecioKg KgUnidad
Set-DuraznoRow 207 44918 "Early Majestic" "Especial" 300 22000 22000 22000 `
    "$/caja 14 kilos empedrada" "Región de O'Higgins" 1571 14

Set-DuraznoRow 208 44918 "Early Majestic" "Primera" 600 18000 19000 18500 `
    "$/caja 14 kilos empedrada" "Región de O'Higgins" 1321 14

Set-DuraznoRow 209 44918 "Florida King" "Especial" 300 22000 22000 22000 `
    "$/caja 14 kilos empedrada" "Región de O'Higgins" 1571 14

Set-DuraznoRow 210 44918 "Florida King" "Primera" 600 18000 19000 18500 `
    "$/caja 14 kilos empedrada" "Región de O'Higgins" 1321 14
